# Auto-generated edit script applying the Louisoix_Profits.xlsx diff
# (values recomputed after a source-data refresh; table formulas are static
# numeric cells, so each changed cell is written directly via the COM object
# model as the diff specifies, including new/removed cells.)
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(57, 8).Value = 58199.2  # H57: 62598.4 -> 58199.2
$ws.Cells.Item(57, 9).Value = 28999  # I57: 0 -> 28999
$ws.Cells.Item(57, 10).Value = 65499.25  # J57: 62598.4 -> 65499.25
$ws.Cells.Item(57, 11).Value = 86997  # K57: 0 -> 86997
$ws.Cells.Item(57, 12).Value = 196497.75  # L57: 187795.2 -> 196497.75
$ws.Cells.Item(57, 13).Value = -86498  # M57: None -> -86498
$ws.Cells.Item(57, 14).Value = -197495.75  # N57: -188793.2 -> -197495.75
$ws.Cells.Item(64, 8).Value = 5019.75  # H64: 5073.1113 -> 5019.75
$ws.Cells.Item(64, 10).Value = 5111.5  # J64: 5167 -> 5111.5
$ws.Cells.Item(64, 12).Value = 5111.5  # L64: 5167 -> 5111.5
$ws.Cells.Item(64, 14).Value = -5607.5  # N64: -5663 -> -5607.5
$ws.Cells.Item(67, 8).Value = 5019.75  # H67: 5073.1113 -> 5019.75
$ws.Cells.Item(67, 10).Value = 5111.5  # J67: 5167 -> 5111.5
$ws.Cells.Item(67, 12).Value = 5111.5  # L67: 5167 -> 5111.5
$ws.Cells.Item(67, 14).Value = -6827.5  # N67: -6883 -> -6827.5
$ws.Cells.Item(86, 8).Value = 2646.3333  # H86: 2914.1875 -> 2646.3333
$ws.Cells.Item(86, 9).Value = 2591.75  # I86: 3172.6667 -> 2591.75
$ws.Cells.Item(86, 10).Value = 2690  # J86: 2759.1 -> 2690
$ws.Cells.Item(86, 11).Value = 2591.75  # K86: 3172.6667 -> 2591.75
$ws.Cells.Item(86, 12).Value = 2690  # L86: 2759.1 -> 2690
$ws.Cells.Item(86, 13).Value = -1468.75  # M86: -2049.6667 -> -1468.75
$ws.Cells.Item(86, 14).Value = -4936  # N86: -5005.1 -> -4936
$ws.Cells.Item(89, 8).Value = 2646.3333  # H89: 2914.1875 -> 2646.3333
$ws.Cells.Item(89, 9).Value = 2591.75  # I89: 3172.6667 -> 2591.75
$ws.Cells.Item(89, 10).Value = 2690  # J89: 2759.1 -> 2690
$ws.Cells.Item(89, 11).Value = 12958.75  # K89: 15863.3335 -> 12958.75
$ws.Cells.Item(89, 12).Value = 13450  # L89: 13795.5 -> 13450
$ws.Cells.Item(89, 13).Value = -7342.75  # M89: -10247.3335 -> -7342.75
$ws.Cells.Item(89, 14).Value = -24682  # N89: -25027.5 -> -24682
$ws.Cells.Item(98, 8).Value = 1785  # H98: 1824.6 -> 1785
$ws.Cells.Item(98, 9).Value = 1191.2424  # I98: 1216 -> 1191.2424
$ws.Cells.Item(98, 11).Value = 1191.2424  # K98: 1216 -> 1191.2424
$ws.Cells.Item(98, 13).Value = 306.7575999999999  # M98: 282 -> 306.7575999999999
$ws.Cells.Item(122, 8).Value = 1785  # H122: 1824.6 -> 1785
$ws.Cells.Item(122, 9).Value = 1191.2424  # I122: 1216 -> 1191.2424
$ws.Cells.Item(122, 11).Value = 3573.7272  # K122: 3648 -> 3573.7272
$ws.Cells.Item(122, 13).Value = -1123.7272  # M122: -1198 -> -1123.7272
$ws.Cells.Item(137, 8).Value = 1462.8948  # H137: 1392.5652 -> 1462.8948
$ws.Cells.Item(137, 9).Value = 1193.8235  # I137: 1133.6316 -> 1193.8235
$ws.Cells.Item(137, 10).Value = 3750  # J137: 2622.5 -> 3750
$ws.Cells.Item(137, 11).Value = 3581.4705  # K137: 3400.8948 -> 3581.4705
$ws.Cells.Item(137, 12).Value = 11250  # L137: 7867.5 -> 11250
$ws.Cells.Item(137, 13).Value = -1031.4705  # M137: -850.8948 -> -1031.4705
$ws.Cells.Item(137, 14).Value = -16350  # N137: -12967.5 -> -16350
$ws.Cells.Item(139, 8).Value = 64466.668  # H139: 72200 -> 64466.668
$ws.Cells.Item(139, 10).Value = 71700  # J139: 94400 -> 71700
$ws.Cells.Item(139, 12).Value = 71700  # L139: 94400 -> 71700
$ws.Cells.Item(139, 14).Value = -81980  # N139: -104680 -> -81980

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2049.4443  # H2: 2054.4443 -> 2049.4443
$ws.Cells.Item(2, 9).Value = 2049.4443  # I2: 2054.4443 -> 2049.4443
$ws.Cells.Item(2, 11).Value = 2049.4443  # K2: 2054.4443 -> 2049.4443
$ws.Cells.Item(2, 13).Value = -1936.4443  # M2: -1941.4443 -> -1936.4443
$ws.Cells.Item(61, 8).Value = 4498.9165  # H61: 4575.1816 -> 4498.9165
$ws.Cells.Item(61, 9).Value = 3923.2222  # I61: 3956.125 -> 3923.2222
$ws.Cells.Item(61, 11).Value = 3923.2222  # K61: 3956.125 -> 3923.2222
$ws.Cells.Item(61, 13).Value = -3711.2222  # M61: -3744.125 -> -3711.2222
$ws.Cells.Item(116, 8).Value = 2049.4443  # H116: 2054.4443 -> 2049.4443
$ws.Cells.Item(116, 9).Value = 2049.4443  # I116: 2054.4443 -> 2049.4443
$ws.Cells.Item(116, 11).Value = 2049.4443  # K116: 2054.4443 -> 2049.4443
$ws.Cells.Item(116, 13).Value = 244.5556999999999  # M116: 239.5556999999999 -> 244.5556999999999
$ws.Cells.Item(136, 8).Value = 4498.9165  # H136: 4575.1816 -> 4498.9165
$ws.Cells.Item(136, 9).Value = 3923.2222  # I136: 3956.125 -> 3923.2222
$ws.Cells.Item(136, 11).Value = 11769.6666  # K136: 11868.375 -> 11769.6666
$ws.Cells.Item(136, 13).Value = -9219.6666  # M136: -9318.375 -> -9219.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2049.4443  # H3: 2054.4443 -> 2049.4443
$ws.Cells.Item(3, 9).Value = 2049.4443  # I3: 2054.4443 -> 2049.4443
$ws.Cells.Item(3, 11).Value = 2049.4443  # K3: 2054.4443 -> 2049.4443
$ws.Cells.Item(3, 13).Value = -1935.4443  # M3: -1940.4443 -> -1935.4443
$ws.Cells.Item(20, 8).Value = 5885340.5  # H20: 6253034 -> 5885340.5
$ws.Cells.Item(20, 9).Value = 12501312  # I20: 12501506 -> 12501312
$ws.Cells.Item(20, 10).Value = 4476.8887  # J20: 4561.5 -> 4476.8887
$ws.Cells.Item(20, 11).Value = 12501312  # K20: 12501506 -> 12501312
$ws.Cells.Item(20, 12).Value = 4476.8887  # L20: 4561.5 -> 4476.8887
$ws.Cells.Item(20, 13).Value = -12501065  # M20: -12501259 -> -12501065
$ws.Cells.Item(20, 14).Value = -4970.8887  # N20: -5055.5 -> -4970.8887
$ws.Cells.Item(26, 8).Value = 27200  # H26: 42104 -> 27200
$ws.Cells.Item(26, 9).Value = 26920.25  # I26: 42632.332 -> 26920.25
$ws.Cells.Item(26, 10).Value = 27759.5  # J26: 40519 -> 27759.5
$ws.Cells.Item(26, 11).Value = 26920.25  # K26: 42632.332 -> 26920.25
$ws.Cells.Item(26, 12).Value = 27759.5  # L26: 40519 -> 27759.5
$ws.Cells.Item(26, 13).Value = -26628.25  # M26: -42340.332 -> -26628.25
$ws.Cells.Item(26, 14).Value = -28343.5  # N26: -41103 -> -28343.5
$ws.Cells.Item(74, 8).Value = 53315.668  # H74: 59981.75 -> 53315.668
$ws.Cells.Item(74, 9).Value = 39974  # I74: 0 -> 39974
$ws.Cells.Item(74, 10).Value = 59986.5  # J74: 59981.75 -> 59986.5
$ws.Cells.Item(74, 11).Value = 39974  # K74: 0 -> 39974
$ws.Cells.Item(74, 12).Value = 59986.5  # L74: 59981.75 -> 59986.5
$ws.Cells.Item(74, 13).Value = -39038  # M74: None -> -39038
$ws.Cells.Item(74, 14).Value = -61858.5  # N74: -61853.75 -> -61858.5
$ws.Cells.Item(77, 8).Value = 53315.668  # H77: 59981.75 -> 53315.668
$ws.Cells.Item(77, 9).Value = 39974  # I77: 0 -> 39974
$ws.Cells.Item(77, 10).Value = 59986.5  # J77: 59981.75 -> 59986.5
$ws.Cells.Item(77, 11).Value = 119922  # K77: 0 -> 119922
$ws.Cells.Item(77, 12).Value = 179959.5  # L77: 179945.25 -> 179959.5
$ws.Cells.Item(77, 13).Value = -115242  # M77: None -> -115242
$ws.Cells.Item(77, 14).Value = -189319.5  # N77: -189305.25 -> -189319.5
$ws.Cells.Item(86, 8).Value = 2849.1904  # H86: 3062 -> 2849.1904
$ws.Cells.Item(86, 9).Value = 1783.3043  # I86: 2038.4 -> 1783.3043
$ws.Cells.Item(86, 11).Value = 1783.3043  # K86: 2038.4 -> 1783.3043
$ws.Cells.Item(86, 13).Value = -660.3043  # M86: -915.4000000000001 -> -660.3043
$ws.Cells.Item(89, 8).Value = 2849.1904  # H89: 3062 -> 2849.1904
$ws.Cells.Item(89, 9).Value = 1783.3043  # I89: 2038.4 -> 1783.3043
$ws.Cells.Item(89, 11).Value = 8916.521500000001  # K89: 10192 -> 8916.521500000001
$ws.Cells.Item(89, 13).Value = -3300.521500000001  # M89: -4576 -> -3300.521500000001
$ws.Cells.Item(94, 8).Value = 1758.2222  # H94: 1725.6786 -> 1758.2222
$ws.Cells.Item(94, 9).Value = 1591.1305  # I94: 1560.125 -> 1591.1305
$ws.Cells.Item(94, 11).Value = 1591.1305  # K94: 1560.125 -> 1591.1305
$ws.Cells.Item(94, 13).Value = -1140.1305  # M94: -1109.125 -> -1140.1305
$ws.Cells.Item(110, 8).Value = 37497.5  # H110: 34997 -> 37497.5
$ws.Cells.Item(110, 10).Value = 37497.5  # J110: 34997 -> 37497.5
$ws.Cells.Item(110, 12).Value = 37497.5  # L110: 34997 -> 37497.5
$ws.Cells.Item(110, 14).Value = -45677.5  # N110: -43177 -> -45677.5
$ws.Cells.Item(134, 8).Value = 3752.147  # H134: 4167.483 -> 3752.147
$ws.Cells.Item(134, 9).Value = 2985.5908  # I134: 3468.647 -> 2985.5908
$ws.Cells.Item(134, 11).Value = 8956.7724  # K134: 10405.941 -> 8956.7724
$ws.Cells.Item(134, 13).Value = -6421.7724  # M134: -7870.940999999999 -> -6421.7724

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 575.93335  # H22: 557.4375 -> 575.93335
$ws.Cells.Item(22, 9).Value = 364.2  # I22: 356.54544 -> 364.2
$ws.Cells.Item(22, 11).Value = 364.2  # K22: 356.54544 -> 364.2
$ws.Cells.Item(22, 13).Value = -14.19999999999999  # M22: -6.545439999999985 -> -14.19999999999999
$ws.Cells.Item(51, 8).Value = 40856.715  # H51: 41333.332 -> 40856.715
$ws.Cells.Item(51, 10).Value = 47199.4  # J51: 49500 -> 47199.4
$ws.Cells.Item(51, 12).Value = 47199.4  # L51: 49500 -> 47199.4
$ws.Cells.Item(51, 14).Value = -48671.4  # N51: -50972 -> -48671.4
$ws.Cells.Item(57, 8).Value = 49998.5  # H57: 49999 -> 49998.5
$ws.Cells.Item(57, 10).Value = 49998.5  # J57: 49999 -> 49998.5
$ws.Cells.Item(57, 12).Value = 49998.5  # L57: 49999 -> 49998.5
$ws.Cells.Item(57, 14).Value = -51118.5  # N57: -51119 -> -51118.5
$ws.Cells.Item(61, 8).Value = 40856.715  # H61: 41333.332 -> 40856.715
$ws.Cells.Item(61, 10).Value = 47199.4  # J61: 49500 -> 47199.4
$ws.Cells.Item(61, 12).Value = 47199.4  # L61: 49500 -> 47199.4
$ws.Cells.Item(61, 14).Value = -47895.4  # N61: -50196 -> -47895.4
$ws.Cells.Item(62, 8).Value = 4258.125  # H62: 4295 -> 4258.125
$ws.Cells.Item(62, 10).Value = 4679  # J62: 4814.8 -> 4679
$ws.Cells.Item(62, 12).Value = 4679  # L62: 4814.8 -> 4679
$ws.Cells.Item(62, 14).Value = -5927  # N62: -6062.8 -> -5927
$ws.Cells.Item(65, 8).Value = 4258.125  # H65: 4295 -> 4258.125
$ws.Cells.Item(65, 10).Value = 4679  # J65: 4814.8 -> 4679
$ws.Cells.Item(65, 12).Value = 23395  # L65: 24074 -> 23395
$ws.Cells.Item(65, 14).Value = -29635  # N65: -30314 -> -29635
$ws.Cells.Item(86, 8).Value = 5750.1665  # H86: 6300.2 -> 5750.1665
$ws.Cells.Item(86, 9).Value = 6175.25  # I86: 7233.6665 -> 6175.25
$ws.Cells.Item(86, 11).Value = 6175.25  # K86: 7233.6665 -> 6175.25
$ws.Cells.Item(86, 13).Value = -5052.25  # M86: -6110.6665 -> -5052.25
$ws.Cells.Item(89, 8).Value = 5750.1665  # H89: 6300.2 -> 5750.1665
$ws.Cells.Item(89, 9).Value = 6175.25  # I89: 7233.6665 -> 6175.25
$ws.Cells.Item(89, 11).Value = 30876.25  # K89: 36168.3325 -> 30876.25
$ws.Cells.Item(89, 13).Value = -25260.25  # M89: -30552.3325 -> -25260.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 13998  # H14: 4696.6665 -> 13998
$ws.Cells.Item(14, 9).Value = 13998  # I14: 4696.6665 -> 13998
$ws.Cells.Item(14, 11).Value = 41994  # K14: 14089.9995 -> 41994
$ws.Cells.Item(14, 13).Value = -41821  # M14: -13916.9995 -> -41821
$ws.Cells.Item(37, 8).Value = 64729.8  # H37: 64877.777 -> 64729.8
$ws.Cells.Item(37, 10).Value = 64729.8  # J37: 64877.777 -> 64729.8
$ws.Cells.Item(37, 12).Value = 194189.4  # L37: 194633.331 -> 194189.4
$ws.Cells.Item(37, 14).Value = -194413.4  # N37: -194857.331 -> -194413.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 39999  # H26: 0 -> 39999
$ws.Cells.Item(26, 10).Value = 39999  # J26: 0 -> 39999
$ws.Cells.Item(26, 12).Value = 39999  # L26: 0 -> 39999
$ws.Cells.Item(26, 14).Value = -40559  # N26: None -> -40559
$ws.Cells.Item(50, 8).Value = 39999  # H50: 0 -> 39999
$ws.Cells.Item(50, 10).Value = 39999  # J50: 0 -> 39999
$ws.Cells.Item(50, 12).Value = 39999  # L50: 0 -> 39999
$ws.Cells.Item(50, 14).Value = -40995  # N50: None -> -40995
$ws.Cells.Item(63, 8).Value = 25057  # H63: 30000 -> 25057
$ws.Cells.Item(63, 10).Value = 25057  # J63: 30000 -> 25057
$ws.Cells.Item(63, 12).Value = 25057  # L63: 30000 -> 25057
$ws.Cells.Item(63, 14).Value = -26429  # N63: -31372 -> -26429
$ws.Cells.Item(66, 8).Value = 25057  # H66: 30000 -> 25057
$ws.Cells.Item(66, 10).Value = 25057  # J66: 30000 -> 25057
$ws.Cells.Item(66, 12).Value = 75171  # L66: 90000 -> 75171
$ws.Cells.Item(66, 14).Value = -82035  # N66: -96864 -> -82035
$ws.Cells.Item(70, 8).Value = 4924.25  # H70: 4932.6665 -> 4924.25
$ws.Cells.Item(70, 10).Value = 4932.6665  # J70: 4949.5 -> 4932.6665
$ws.Cells.Item(70, 12).Value = 4932.6665  # L70: 4949.5 -> 4932.6665
$ws.Cells.Item(70, 14).Value = -5472.6665  # N70: -5489.5 -> -5472.6665
$ws.Cells.Item(73, 8).Value = 4924.25  # H73: 4932.6665 -> 4924.25
$ws.Cells.Item(73, 10).Value = 4932.6665  # J73: 4949.5 -> 4932.6665
$ws.Cells.Item(73, 12).Value = 4932.6665  # L73: 4949.5 -> 4932.6665
$ws.Cells.Item(73, 14).Value = -6804.6665  # N73: -6821.5 -> -6804.6665
$ws.Cells.Item(113, 8).Value = 160304.69  # H113: 122805.82 -> 160304.69
$ws.Cells.Item(113, 9).Value = 117049.78  # I113: 88023.914 -> 117049.78
$ws.Cells.Item(113, 10).Value = 257628.25  # J113: 206282.4 -> 257628.25
$ws.Cells.Item(113, 11).Value = 117049.78  # K113: 88023.914 -> 117049.78
$ws.Cells.Item(113, 12).Value = 257628.25  # L113: 206282.4 -> 257628.25
$ws.Cells.Item(113, 13).Value = -114879.78  # M113: -85853.914 -> -114879.78
$ws.Cells.Item(113, 14).Value = -261968.25  # N113: -210622.4 -> -261968.25
$ws.Cells.Item(132, 8).Value = 49420.715  # H132: 47287.863 -> 49420.715
$ws.Cells.Item(132, 9).Value = 54228.26  # I132: 57185.723 -> 54228.26
$ws.Cells.Item(132, 10).Value = 3749  # J132: 2747.5 -> 3749
$ws.Cells.Item(132, 11).Value = 162684.78  # K132: 171557.169 -> 162684.78
$ws.Cells.Item(132, 12).Value = 11247  # L132: 8242.5 -> 11247
$ws.Cells.Item(132, 13).Value = -160154.78  # M132: -169027.169 -> -160154.78
$ws.Cells.Item(132, 14).Value = -16307  # N132: -13302.5 -> -16307

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5261.625  # H7: 4253.5 -> 5261.625
$ws.Cells.Item(7, 9).Value = 2418.8  # I7: 2338.111 -> 2418.8
$ws.Cells.Item(7, 11).Value = 2418.8  # K7: 2338.111 -> 2418.8
$ws.Cells.Item(7, 13).Value = -2306.8  # M7: -2226.111 -> -2306.8
$ws.Cells.Item(46, 8).Value = 2062.9395  # H46: 2338.3928 -> 2062.9395
$ws.Cells.Item(46, 9).Value = 975.5  # I46: 1500 -> 975.5
$ws.Cells.Item(46, 10).Value = 2212.9312  # J46: 2369.4443 -> 2212.9312
$ws.Cells.Item(46, 11).Value = 975.5  # K46: 1500 -> 975.5
$ws.Cells.Item(46, 12).Value = 2212.9312  # L46: 2369.4443 -> 2212.9312
$ws.Cells.Item(46, 13).Value = -787.5  # M46: -1312 -> -787.5
$ws.Cells.Item(46, 14).Value = -2588.9312  # N46: -2745.4443 -> -2588.9312
$ws.Cells.Item(82, 8).Value = 2218.92  # H82: 2211.923 -> 2218.92
$ws.Cells.Item(82, 9).Value = 2014.7693  # I82: 2031.2858 -> 2014.7693
$ws.Cells.Item(82, 10).Value = 2440.0833  # J82: 2422.6667 -> 2440.0833
$ws.Cells.Item(82, 11).Value = 2014.7693  # K82: 2031.2858 -> 2014.7693
$ws.Cells.Item(82, 12).Value = 2440.0833  # L82: 2422.6667 -> 2440.0833
$ws.Cells.Item(82, 13).Value = -1653.7693  # M82: -1670.2858 -> -1653.7693
$ws.Cells.Item(82, 14).Value = -3162.0833  # N82: -3144.6667 -> -3162.0833
$ws.Cells.Item(85, 8).Value = 2218.92  # H85: 2211.923 -> 2218.92
$ws.Cells.Item(85, 9).Value = 2014.7693  # I85: 2031.2858 -> 2014.7693
$ws.Cells.Item(85, 10).Value = 2440.0833  # J85: 2422.6667 -> 2440.0833
$ws.Cells.Item(85, 11).Value = 2014.7693  # K85: 2031.2858 -> 2014.7693
$ws.Cells.Item(85, 12).Value = 2440.0833  # L85: 2422.6667 -> 2440.0833
$ws.Cells.Item(85, 13).Value = -766.7692999999999  # M85: -783.2858000000001 -> -766.7692999999999
$ws.Cells.Item(85, 14).Value = -4936.0833  # N85: -4918.6667 -> -4936.0833
$ws.Cells.Item(110, 8).Value = 0  # H110: 644 -> 0
$ws.Cells.Item(110, 10).Value = 0  # J110: 644 -> 0
$ws.Cells.Item(110, 12).Value = 0  # L110: 644 -> 0
$ws.Cells.Item(110, 14).ClearContents()  # N110: remove (was -8824)
$ws.Cells.Item(126, 8).Value = 5261.625  # H126: 4253.5 -> 5261.625
$ws.Cells.Item(126, 9).Value = 2418.8  # I126: 2338.111 -> 2418.8
$ws.Cells.Item(126, 11).Value = 7256.400000000001  # K126: 7014.333 -> 7256.400000000001
$ws.Cells.Item(126, 13).Value = -4786.400000000001  # M126: -4544.333 -> -4786.400000000001
$ws.Cells.Item(136, 8).Value = 6417.2  # H136: 5818.5 -> 6417.2
$ws.Cells.Item(136, 9).Value = 1089.5  # I136: 1059.6666 -> 1089.5
$ws.Cells.Item(136, 10).Value = 7749.125  # J136: 7858 -> 7749.125
$ws.Cells.Item(136, 11).Value = 3268.5  # K136: 3178.9998 -> 3268.5
$ws.Cells.Item(136, 12).Value = 23247.375  # L136: 23574 -> 23247.375
$ws.Cells.Item(136, 13).Value = -718.5  # M136: -628.9998000000001 -> -718.5
$ws.Cells.Item(136, 14).Value = -28347.375  # N136: -28674 -> -28347.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 47450.137  # H132: 64900.938 -> 47450.137
$ws.Cells.Item(132, 9).Value = 73459.28999999999  # I132: 113759.11 -> 73459.28999999999
$ws.Cells.Item(132, 10).Value = 1934.125  # J132: 2083.2856 -> 1934.125
$ws.Cells.Item(132, 11).Value = 220377.87  # K132: 341277.33 -> 220377.87
$ws.Cells.Item(132, 12).Value = 5802.375  # L132: 6249.8568 -> 5802.375
$ws.Cells.Item(132, 13).Value = -217847.87  # M132: -338747.33 -> -217847.87
$ws.Cells.Item(132, 14).Value = -10862.375  # N132: -11309.8568 -> -10862.375
$ws.Cells.Item(133, 8).Value = 89999  # H133: 90000 -> 89999
$ws.Cells.Item(133, 10).Value = 89999  # J133: 90000 -> 89999
$ws.Cells.Item(133, 12).Value = 89999  # L133: 90000 -> 89999
$ws.Cells.Item(133, 14).Value = -100119  # N133: -100120 -> -100119
